$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 104, shifting rows 104:113 down to 105:114
$ws.Rows.Item(104).Insert()

# Populate the newly inserted row 104 with data (copy static columns from row 105,
# which now holds what used to be in row 104, and set the specific new values).
$ws.Range("A104").Value = 4
$ws.Range("B104").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C104").Value = "Los Lagos"
$ws.Range("D104").Value = 44826
$ws.Range("E104").Value = 10
$ws.Range("F104").Value = "Fruta"
$ws.Range("G104").Value = 100104
$ws.Range("H104").Value = "Frutos de pepita"
$ws.Range("I104").Value = 100104003
$ws.Range("J104").Value = "Membrillo"
$ws.Range("K104").Value = "Champion"
$ws.Range("L104").Value = "Primera"
$ws.Range("M104").Value = 60
$ws.Range("N104").Value = 14000
$ws.Range("O104").Value = 15000
$ws.Range("P104").Value = 14500
$ws.Range("Q104").Value = "`$/caja 18 kilos granel"
$ws.Range("R104").Value = "Región de O'Higgins"
$ws.Range("S104").Value = 806
$ws.Range("T104").Value = 18
